$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values added to column C (and E) for rows 4-6, plus updated values
# for rows 4-19 across columns C and E.
$ws.Range("C4").Value = 0.7825601129312298
$ws.Range("E4").Value = -0.2328395085068102

$ws.Range("C5").Value = 1.171373351779592
$ws.Range("E5").Value = 0.6458049593451864

$ws.Range("C6").Value = 0.4797371259343874
$ws.Range("E6").Value = 0.5448268972618964

$ws.Range("C7").Value = 0.1088602047940146
$ws.Range("E7").Value = 0.1671551101610103

$ws.Range("C8").Value = -0.075394216261504
$ws.Range("E8").Value = -0.06463514052835739

$ws.Range("C9").Value = -0.07548837955325682
$ws.Range("E9").Value = 0.05928147027902675

$ws.Range("C10").Value = -0.5735475396625112
$ws.Range("E10").Value = 0.02406984837131088

$ws.Range("C11").Value = 0.02883110668334687
$ws.Range("E11").Value = 0.241498802789164

$ws.Range("C12").Value = 0.7252300059688022
$ws.Range("E12").Value = 0.2052430644269299

$ws.Range("C13").Value = -0.6243248145489155
$ws.Range("E13").Value = 0.07475225043114264

$ws.Range("C14").Value = -0.6491730431770759
$ws.Range("E14").Value = -0.4370777949570193

$ws.Range("C15").Value = 1.311904119834839
$ws.Range("E15").Value = -0.1159018519404809

$ws.Range("C16").Value = -1.338216592160768
$ws.Range("E16").Value = -0.2487719682984557

$ws.Range("C17").Value = 0.2512652100014945
$ws.Range("E17").Value = 0.041441321352087

$ws.Range("C18").Value = 0.9693451788297391
$ws.Range("E18").Value = 0.08029846083614789

$ws.Range("C19").Value = -1.551451534890558
$ws.Range("E19").Value = -0.1079783528070921
